$wb = $excel.ActiveWorkbook

# Sheet "建物" (Building): column I row 2 holds the "property_category"
# value, which was incorrectly set to "land". Fix it to "building".
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"

# Sheet "汽車" (Car): column H row 2 holds the "property_category"
# value, which was incorrectly set to "land". Fix it to "car".
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
